# Update the lattice-multiplication exercise table: each cell keeps its
# layout (problem header, split digits, divider, two partial-product rows)
# but the numbers themselves change to a new set of exercises.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

function Set-CellText($row, $col, $lines) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = ($lines -join $nl)
}

Set-CellText 1 1 @("84 x 90", "  9    0", "  ----", "8|    |", "4|    |")
Set-CellText 1 2 @("92 x 32", "  3    2", "  ----", "9|    |", "2|    |")
Set-CellText 1 3 @("19 x 52", "  5    2", "  ----", "1|    |", "9|    |")

Set-CellText 2 1 @("48 x 71", "  7    1", "  ----", "4|    |", "8|    |")
Set-CellText 2 2 @("11 x 95", "  9    5", "  ----", "1|    |", "1|    |")
Set-CellText 2 3 @("94 x 24", "  2    4", "  ----", "9|    |", "4|    |")

Set-CellText 3 1 @("28 x 15", "  1    5", "  ----", "2|    |", "8|    |")
Set-CellText 3 2 @("45 x 72", "  7    2", "  ----", "4|    |", "5|    |")
Set-CellText 3 3 @("59 x 38", "  3    8", "  ----", "5|    |", "9|    |")

Set-CellText 4 1 @("59 x 55", "  5    5", "  ----", "5|    |", "9|    |")
Set-CellText 4 2 @("56 x 83", "  8    3", "  ----", "5|    |", "6|    |")
Set-CellText 4 3 @("41 x 21", "  2    1", "  ----", "4|    |", "1|    |")

Set-CellText 5 1 @("54 x 87", "  8    7", "  ----", "5|    |", "4|    |")
Set-CellText 5 2 @("69 x 67", "  6    7", "  ----", "6|    |", "9|    |")
Set-CellText 5 3 @("76 x 75", "  7    5", "  ----", "7|    |", "6|    |")

Write-Output "lattice multiplication exercises updated"
